$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List2")

# 1) Add new AQ (raw data) / AR (GEH formula) columns for rows 2-25
$values = 126,100,88,45,49,114,183,322,528,841,1177,1314,1422,1625,1793,1908,1779,1769,1793,1800,1028,580,350,237
$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 43).Value = $v
    $ws.Cells.Item($row, 44).Formula = "=SQRT(2*(AQ" + $row + "-B" + $row + ")^2/(AQ" + $row + "+B" + $row + "))"
    $row++
}
# Row 25's GEH formula column (AR) keeps no border/style (matches neighbouring result columns K,L,M,N / AC..AF / AL..AO)
$ws.Cells.Item(25, 44).Borders.LineStyle = 0
# AS25 exists as an empty, unstyled cell (extends used range to column AS)
$ws.Range("AS25").Borders.LineStyle = 0

# 2) Apply conditional formatting to the new AS (blank) and AR (GEH) columns, mirroring the existing scheme
$rngAS = $ws.Range("AS2:AS25")
$fcAS1 = $rngAS.FormatConditions.Add(1, 5, "10")
$fcAS1.Font.Color = 393372
$fcAS1.Interior.Color = 13551615
$fcAS2 = $rngAS.FormatConditions.AddColorScale(3)
$fcAS3 = $rngAS.FormatConditions.Add(1, 5, "10")

$rngAR = $ws.Range("AR2:AR25")
$fcAR1 = $rngAR.FormatConditions.Add(1, 5, "10")
$fcAR1.Font.Color = 393372
$fcAR1.Interior.Color = 13551615
$fcAR2 = $rngAR.FormatConditions.AddColorScale(3)
$fcAR3 = $rngAR.FormatConditions.Add(1, 5, "10")

# 3) Renumber priorities: new AR block -> 1,2,3 ; new AS block -> 4,5,6 ; existing 7 blocks shift up by 6
$fcAR1.Priority = 1
$fcAR2.Priority = 2
$fcAR3.Priority = 3
$fcAS1.Priority = 4
$fcAS2.Priority = 5
$fcAS3.Priority = 6

$existing = @(
    @{addr="K26:K1048576"; base=19},
    @{addr="T2:W25"; base=16},
    @{addr="AC2:AF25"; base=13},
    @{addr="T26:W26"; base=10},
    @{addr="AC26:AF26"; base=7},
    @{addr="K2:N25"; base=4},
    @{addr="AL2:AO25"; base=1}
)
foreach ($item in $existing) {
    $r = $ws.Range($item.addr)
    $fcs = $r.FormatConditions
    for ($i=1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).Priority = $item.base + ($i-1) + 6
    }
}

# 4) Update the view: scroll position and selected cell
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 26
$ws.Range("AR21").Select()
